# "stubble sim" land-use mix table: add a mixed-species land-use group.
# This inserts a new 3-column group (value / weight / value) after the
# existing "AD" group, shifting every later group right by 3 columns,
# and recalculates the weighting for the group now at columns V:X.
# The new group is currently turned off (weight ~= the "off" sentinel),
# still needs calibration, and causes no profit change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object "object[,]" 10,54

# Row 1
$values[0,0] = 0.04887908386567361
$values[0,1] = 0.0001146235521235521
$values[0,2] = 0.04887908386567361
$values[0,3] = 0.04887908386567361
$values[0,4] = 0.0001146235521235521
$values[0,5] = 0.04887908386567361
$values[0,6] = 0.03528994484009806
$values[0,7] = 0.0001518968383129674
$values[0,8] = 0.03528994484009806
$values[0,9] = 0.03550473685604939
$values[0,10] = 0.0001146235521235521
$values[0,11] = 0.03550473685604939
$values[0,12] = 0.03528994484009806
$values[0,13] = 0.0001518968383129674
$values[0,14] = 0.03528994484009806
$values[0,15] = 0.03528994484009806
$values[0,16] = 0.0001518968383129674
$values[0,17] = 0.03528994484009806
$values[0,18] = 0.03528994484009806
$values[0,19] = 0.0001518968383129674
$values[0,20] = 0.03528994484009806
$values[0,21] = 0.03003653031203316
$values[0,22] = 0.0001394724095831623
$values[0,23] = 0.03003653031203316
$values[0,24] = 0.03550473685604939
$values[0,25] = 0.0001146235521235521
$values[0,26] = 0.03550473685604939
$values[0,27] = 0.03550473685604939
$values[0,28] = 0.0001146235521235521
$values[0,29] = 0.03550473685604939
$values[0,30] = 0.03550473685604939
$values[0,31] = 0.0001146235521235521
$values[0,32] = 0.03550473685604939
$values[0,33] = 0.01647340995996823
$values[0,34] = 0.0001146235521235521
$values[0,35] = 0.01647340995996823
$values[0,36] = 0.01647340995996823
$values[0,37] = 0.0001146235521235521
$values[0,38] = 0.01647340995996823
$values[0,39] = 0.03528994484009806
$values[0,40] = 0.0001518968383129674
$values[0,41] = 0.03528994484009806
$values[0,42] = 0.03104662118617465
$values[0,43] = 0.0001146235521235521
$values[0,44] = 0.03104662118617465
$values[0,45] = 0.03104662118617465
$values[0,46] = 0.0001146235521235521
$values[0,47] = 0.03104662118617465
$values[0,48] = 0.01647340995996823
$values[0,49] = 0.0001146235521235521
$values[0,50] = 0.01647340995996823
$values[0,51] = 0.01647340995996823
$values[0,52] = 0.0001146235521235521
$values[0,53] = 0.01647340995996823

# Row 2
$values[1,0] = 0.006729471584568375
$values[1,1] = 0.0001146235521235521
$values[1,2] = 0.006729471584568375
$values[1,3] = 0.006729471584568375
$values[1,4] = 0.0001146235521235521
$values[1,5] = 0.006729471584568375
$values[1,6] = 0.008613748025750392
$values[1,7] = 0.1726507920312587
$values[1,8] = 0.008613748025750392
$values[1,9] = 0.007971046728872989
$values[1,10] = 0.0001146235521235521
$values[1,11] = 0.007971046728872989
$values[1,12] = 0.008613748025750392
$values[1,13] = 0.1726507920312587
$values[1,14] = 0.008613748025750392
$values[1,15] = 0.008613748025750392
$values[1,16] = 0.1726507920312587
$values[1,17] = 0.008613748025750392
$values[1,18] = 0.008613748025750392
$values[1,19] = 0.1726507920312587
$values[1,20] = 0.008613748025750392
$values[1,21] = 0.007800672337179483
$values[1,22] = 0.115138735871547
$values[1,23] = 0.007800672337179483
$values[1,24] = 0.007971046728872989
$values[1,25] = 0.0001146235521235521
$values[1,26] = 0.007971046728872989
$values[1,27] = 0.007971046728872989
$values[1,28] = 0.0001146235521235521
$values[1,29] = 0.007971046728872989
$values[1,30] = 0.007971046728872989
$values[1,31] = 0.0001146235521235521
$values[1,32] = 0.007971046728872989
$values[1,33] = 0.007053385709507021
$values[1,34] = 0.0001146235521235521
$values[1,35] = 0.007053385709507021
$values[1,36] = 0.007053385709507021
$values[1,37] = 0.0001146235521235521
$values[1,38] = 0.007053385709507021
$values[1,39] = 0.008613748025750392
$values[1,40] = 0.1726507920312587
$values[1,41] = 0.008613748025750392
$values[1,42] = 0.008384905110307858
$values[1,43] = 0.0001146235521235521
$values[1,44] = 0.008384905110307858
$values[1,45] = 0.008384905110307858
$values[1,46] = 0.0001146235521235521
$values[1,47] = 0.008384905110307858
$values[1,48] = 0.007053385709507021
$values[1,49] = 0.0001146235521235521
$values[1,50] = 0.007053385709507021
$values[1,51] = 0.007053385709507021
$values[1,52] = 0.0001146235521235521
$values[1,53] = 0.007053385709507021

# Row 3
$values[2,0] = 0.0060372868083389
$values[2,1] = 0.3402453775823806
$values[2,2] = 0.0060372868083389
$values[2,3] = 0.0060372868083389
$values[2,4] = 0.3402453775823806
$values[2,5] = 0.0060372868083389
$values[2,6] = 0.005854187974434026
$values[2,7] = 0.3338558167167581
$values[2,8] = 0.005854187974434026
$values[2,9] = 0.006442998431046862
$values[2,10] = 0.3402453775823806
$values[2,11] = 0.006442998431046862
$values[2,12] = 0.005854187974434026
$values[2,13] = 0.3338558167167581
$values[2,14] = 0.005854187974434026
$values[2,15] = 0.005854187974434026
$values[2,16] = 0.3338558167167581
$values[2,17] = 0.005854187974434026
$values[2,18] = 0.005854187974434026
$values[2,19] = 0.3338558167167581
$values[2,20] = 0.005854187974434026
$values[2,21] = 0.005738715325018062
$values[2,22] = 0.3359856703386322
$values[2,23] = 0.005738715325018062
$values[2,24] = 0.006442998431046862
$values[2,25] = 0.3402453775823806
$values[2,26] = 0.006442998431046862
$values[2,27] = 0.006442998431046862
$values[2,28] = 0.3402453775823806
$values[2,29] = 0.006442998431046862
$values[2,30] = 0.006442998431046862
$values[2,31] = 0.3402453775823806
$values[2,32] = 0.006442998431046862
$values[2,33] = 0.004272412788252482
$values[2,34] = 0.3402453775823806
$values[2,35] = 0.004272412788252482
$values[2,36] = 0.004272412788252482
$values[2,37] = 0.3402453775823806
$values[2,38] = 0.004272412788252482
$values[2,39] = 0.005854187974434026
$values[2,40] = 0.3338558167167581
$values[2,41] = 0.005854187974434026
$values[2,42] = 0.006578235638616182
$values[2,43] = 0.3402453775823806
$values[2,44] = 0.006578235638616182
$values[2,45] = 0.006578235638616182
$values[2,46] = 0.3402453775823806
$values[2,47] = 0.006578235638616182
$values[2,48] = 0.004272412788252482
$values[2,49] = 0.3402453775823806
$values[2,50] = 0.004272412788252482
$values[2,51] = 0.004272412788252482
$values[2,52] = 0.3402453775823806
$values[2,53] = 0.004272412788252482

# Row 4
$values[3,0] = 0.005413263642548172
$values[3,1] = 0.2180215461683213
$values[3,2] = 0.005413263642548172
$values[3,3] = 0.005413263642548172
$values[3,4] = 0.2180215461683213
$values[3,5] = 0.005413263642548172
$values[3,6] = 0.003498571737094966
$values[3,7] = 0.1648490742215563
$values[3,8] = 0.003498571737094966
$values[3,9] = 0.004180966898970304
$values[3,10] = 0.2180215461683213
$values[3,11] = 0.004180966898970304
$values[3,12] = 0.003498571737094966
$values[3,13] = 0.1648490742215563
$values[3,14] = 0.003498571737094966
$values[3,15] = 0.003498571737094966
$values[3,16] = 0.1648490742215563
$values[3,17] = 0.003498571737094966
$values[3,18] = 0.003498571737094966
$values[3,19] = 0.1648490742215563
$values[3,20] = 0.003498571737094966
$values[3,21] = 0.003797826962670339
$values[3,22] = 0.1825732315371447
$values[3,23] = 0.003797826962670339
$values[3,24] = 0.004180966898970304
$values[3,25] = 0.2180215461683213
$values[3,26] = 0.004180966898970304
$values[3,27] = 0.004180966898970304
$values[3,28] = 0.2180215461683213
$values[3,29] = 0.004180966898970304
$values[3,30] = 0.004180966898970304
$values[3,31] = 0.2180215461683213
$values[3,32] = 0.004180966898970304
$values[3,33] = 0.003181174702858095
$values[3,34] = 0.2180215461683213
$values[3,35] = 0.003181174702858095
$values[3,36] = 0.003181174702858095
$values[3,37] = 0.2180215461683213
$values[3,38] = 0.003181174702858095
$values[3,39] = 0.003498571737094966
$values[3,40] = 0.1648490742215563
$values[3,41] = 0.003498571737094966
$values[3,42] = 0.003770201317777681
$values[3,43] = 0.2180215461683213
$values[3,44] = 0.003770201317777681
$values[3,45] = 0.003770201317777681
$values[3,46] = 0.2180215461683213
$values[3,47] = 0.003770201317777681
$values[3,48] = 0.003181174702858095
$values[3,49] = 0.2180215461683213
$values[3,50] = 0.003181174702858095
$values[3,51] = 0.003181174702858095
$values[3,52] = 0.2180215461683213
$values[3,53] = 0.003181174702858095

# Row 5
$values[4,0] = 0.004789240476757439
$values[4,1] = 0.1548463323034465
$values[4,2] = 0.004789240476757439
$values[4,3] = 0.004789240476757439
$values[4,4] = 0.1548463323034465
$values[4,5] = 0.004789240476757439
$values[4,6] = 0.004636429019548354
$values[4,7] = 0.1205545666743273
$values[4,8] = 0.004636429019548354
$values[4,9] = 0.004585937853419563
$values[4,10] = 0.1548463323034465
$values[4,11] = 0.004585937853419563
$values[4,12] = 0.004636429019548354
$values[4,13] = 0.1205545666743273
$values[4,14] = 0.004636429019548354
$values[4,15] = 0.004636429019548354
$values[4,16] = 0.1205545666743273
$values[4,17] = 0.004636429019548354
$values[4,18] = 0.004636429019548354
$values[4,19] = 0.1205545666743273
$values[4,20] = 0.004636429019548354
$values[4,21] = 0.00424410917389367
$values[4,22] = 0.131985155217367
$values[4,23] = 0.00424410917389367
$values[4,24] = 0.004585937853419563
$values[4,25] = 0.1548463323034465
$values[4,26] = 0.004585937853419563
$values[4,27] = 0.004585937853419563
$values[4,28] = 0.1548463323034465
$values[4,29] = 0.004585937853419563
$values[4,30] = 0.004585937853419563
$values[4,31] = 0.1548463323034465
$values[4,32] = 0.004585937853419563
$values[4,33] = 0.00336429189201454
$values[4,34] = 0.1548463323034465
$values[4,35] = 0.00336429189201454
$values[4,36] = 0.00336429189201454
$values[4,37] = 0.1548463323034465
$values[4,38] = 0.00336429189201454
$values[4,39] = 0.004636429019548354
$values[4,40] = 0.1205545666743273
$values[4,41] = 0.004636429019548354
$values[4,42] = 0.004518170312306937
$values[4,43] = 0.1548463323034465
$values[4,44] = 0.004518170312306937
$values[4,45] = 0.004518170312306937
$values[4,46] = 0.1548463323034465
$values[4,47] = 0.004518170312306937
$values[4,48] = 0.00336429189201454
$values[4,49] = 0.1548463323034465
$values[4,50] = 0.00336429189201454
$values[4,51] = 0.00336429189201454
$values[4,52] = 0.1548463323034465
$values[4,53] = 0.00336429189201454

# Row 6
$values[5,0] = 0.004119109618312391
$values[5,1] = 0.07167803657229711
$values[5,2] = 0.004119109618312391
$values[5,3] = 0.004119109618312391
$values[5,4] = 0.07167803657229711
$values[5,5] = 0.004119109618312391
$values[5,6] = 0.004031134973236498
$values[5,7] = 0.1485869724678432
$values[5,8] = 0.004031134973236498
$values[5,9] = 0.004416967892344595
$values[5,10] = 0.07167803657229711
$values[5,11] = 0.004416967892344595
$values[5,12] = 0.004031134973236498
$values[5,13] = 0.1485869724678432
$values[5,14] = 0.004031134973236498
$values[5,15] = 0.004031134973236498
$values[5,16] = 0.1485869724678432
$values[5,17] = 0.004031134973236498
$values[5,18] = 0.004031134973236498
$values[5,19] = 0.1485869724678432
$values[5,20] = 0.004031134973236498
$values[5,21] = 0.003936186909500405
$values[5,22] = 0.1229506605026612
$values[5,23] = 0.003936186909500405
$values[5,24] = 0.004416967892344595
$values[5,25] = 0.07167803657229711
$values[5,26] = 0.004416967892344595
$values[5,27] = 0.004416967892344595
$values[5,28] = 0.07167803657229711
$values[5,29] = 0.004416967892344595
$values[5,30] = 0.004416967892344595
$values[5,31] = 0.07167803657229711
$values[5,32] = 0.004416967892344595
$values[5,33] = 0.00292715091194398
$values[5,34] = 0.07167803657229711
$values[5,35] = 0.00292715091194398
$values[5,36] = 0.00292715091194398
$values[5,37] = 0.07167803657229711
$values[5,38] = 0.00292715091194398
$values[5,39] = 0.004031134973236498
$values[5,40] = 0.1485869724678432
$values[5,41] = 0.004031134973236498
$values[5,42] = 0.004516253983688661
$values[5,43] = 0.07167803657229711
$values[5,44] = 0.004516253983688661
$values[5,45] = 0.004516253983688661
$values[5,46] = 0.07167803657229711
$values[5,47] = 0.004516253983688661
$values[5,48] = 0.00292715091194398
$values[5,49] = 0.07167803657229711
$values[5,50] = 0.00292715091194398
$values[5,51] = 0.00292715091194398
$values[5,52] = 0.07167803657229711
$values[5,53] = 0.00292715091194398

# Row 7
$values[6,0] = 0.00349917072984269
$values[6,1] = 0.0001146235521235521
$values[6,2] = 0.00349917072984269
$values[6,3] = 0.00349917072984269
$values[6,4] = 0.0001146235521235521
$values[6,5] = 0.00349917072984269
$values[6,6] = 0.00341847216176009
$values[6,7] = 0.03706360975065923
$values[6,8] = 0.00341847216176009
$values[6,9] = 0.004177719588623684
$values[6,10] = 0.0001146235521235521
$values[6,11] = 0.004177719588623684
$values[6,12] = 0.00341847216176009
$values[6,13] = 0.03706360975065923
$values[6,14] = 0.00341847216176009
$values[6,15] = 0.00341847216176009
$values[6,16] = 0.03706360975065923
$values[6,17] = 0.00341847216176009
$values[6,18] = 0.00341847216176009
$values[6,19] = 0.03706360975065923
$values[6,20] = 0.00341847216176009
$values[6,21] = 0.003708715660526274
$values[6,22] = 0.024747281017814
$values[6,23] = 0.003708715660526274
$values[6,24] = 0.004177719588623684
$values[6,25] = 0.0001146235521235521
$values[6,26] = 0.004177719588623684
$values[6,27] = 0.004177719588623684
$values[6,28] = 0.0001146235521235521
$values[6,29] = 0.004177719588623684
$values[6,30] = 0.004177719588623684
$values[6,31] = 0.0001146235521235521
$values[6,32] = 0.004177719588623684
$values[6,33] = 0.002915829553714546
$values[6,34] = 0.0001146235521235521
$values[6,35] = 0.002915829553714546
$values[6,36] = 0.002915829553714546
$values[6,37] = 0.0001146235521235521
$values[6,38] = 0.002915829553714546
$values[6,39] = 0.00341847216176009
$values[6,40] = 0.03706360975065923
$values[6,41] = 0.00341847216176009
$values[6,42] = 0.004403902541550681
$values[6,43] = 0.0001146235521235521
$values[6,44] = 0.004403902541550681
$values[6,45] = 0.004403902541550681
$values[6,46] = 0.0001146235521235521
$values[6,47] = 0.004403902541550681
$values[6,48] = 0.002915829553714546
$values[6,49] = 0.0001146235521235521
$values[6,50] = 0.002915829553714546
$values[6,51] = 0.002915829553714546
$values[6,52] = 0.0001146235521235521
$values[6,53] = 0.002915829553714546

# Row 8
$values[7,0] = 0.001439616018507967
$values[7,1] = 0.0001146235521235521
$values[7,2] = 0.001439616018507967
$values[7,3] = 0.001439616018507967
$values[7,4] = 0.0001146235521235521
$values[7,5] = 0.001439616018507967
$values[7,6] = 0.002802853103701819
$values[7,7] = 0.02198347762265828
$values[7,8] = 0.002802853103701819
$values[7,9] = 0.002217044442022792
$values[7,10] = 0.0001146235521235521
$values[7,11] = 0.002217044442022792
$values[7,12] = 0.002802853103701819
$values[7,13] = 0.02198347762265828
$values[7,14] = 0.002802853103701819
$values[7,15] = 0.002802853103701819
$values[7,16] = 0.02198347762265828
$values[7,17] = 0.002802853103701819
$values[7,18] = 0.002802853103701819
$values[7,19] = 0.02198347762265828
$values[7,20] = 0.002802853103701819
$values[7,21] = 0.002159886817805626
$values[7,22] = 0.01469385959914671
$values[7,23] = 0.002159886817805626
$values[7,24] = 0.002217044442022792
$values[7,25] = 0.0001146235521235521
$values[7,26] = 0.002217044442022792
$values[7,27] = 0.002217044442022792
$values[7,28] = 0.0001146235521235521
$values[7,29] = 0.002217044442022792
$values[7,30] = 0.002217044442022792
$values[7,31] = 0.0001146235521235521
$values[7,32] = 0.002217044442022792
$values[7,33] = 0.001724088426423199
$values[7,34] = 0.0001146235521235521
$values[7,35] = 0.001724088426423199
$values[7,36] = 0.001724088426423199
$values[7,37] = 0.0001146235521235521
$values[7,38] = 0.001724088426423199
$values[7,39] = 0.002802853103701819
$values[7,40] = 0.02198347762265828
$values[7,41] = 0.002802853103701819
$values[7,42] = 0.002476187249861067
$values[7,43] = 0.0001146235521235521
$values[7,44] = 0.002476187249861067
$values[7,45] = 0.002476187249861067
$values[7,46] = 0.0001146235521235521
$values[7,47] = 0.002476187249861067
$values[7,48] = 0.001724088426423199
$values[7,49] = 0.0001146235521235521
$values[7,50] = 0.001724088426423199
$values[7,51] = 0.001724088426423199
$values[7,52] = 0.0001146235521235521
$values[7,53] = 0.001724088426423199

# Row 9
$values[8,0] = 0.001129646574273117
$values[8,1] = 0.0001146235521235521
$values[8,2] = 0.001129646574273117
$values[8,3] = 0.001129646574273117
$values[8,4] = 0.0001146235521235521
$values[8,5] = 0.001129646574273117
$values[8,6] = 0.0007276711044972345
$values[8,7] = 0.0001518968383129674
$values[8,8] = 0.0007276711044972345
$values[8,9] = 0.001208126731188766
$values[8,10] = 0.0001146235521235521
$values[8,11] = 0.001208126731188766
$values[8,12] = 0.0007276711044972345
$values[8,13] = 0.0001518968383129674
$values[8,14] = 0.0007276711044972345
$values[8,15] = 0.0007276711044972345
$values[8,16] = 0.0001518968383129674
$values[8,17] = 0.0007276711044972345
$values[8,18] = 0.0007276711044972345
$values[8,19] = 0.0001518968383129674
$values[8,20] = 0.0007276711044972345
$values[8,21] = 0.00102508978087982
$values[8,22] = 0.0001394724095831623
$values[8,23] = 0.00102508978087982
$values[8,24] = 0.001208126731188766
$values[8,25] = 0.0001146235521235521
$values[8,26] = 0.001208126731188766
$values[8,27] = 0.001208126731188766
$values[8,28] = 0.0001146235521235521
$values[8,29] = 0.001208126731188766
$values[8,30] = 0.001208126731188766
$values[8,31] = 0.0001146235521235521
$values[8,32] = 0.001208126731188766
$values[8,33] = 0.0008077252184532221
$values[8,34] = 0.0001146235521235521
$values[8,35] = 0.0008077252184532221
$values[8,36] = 0.0008077252184532221
$values[8,37] = 0.0001146235521235521
$values[8,38] = 0.0008077252184532221
$values[8,39] = 0.0007276711044972345
$values[8,40] = 0.0001518968383129674
$values[8,41] = 0.0007276711044972345
$values[8,42] = 0.001234286783493983
$values[8,43] = 0.0001146235521235521
$values[8,44] = 0.001234286783493983
$values[8,45] = 0.001234286783493983
$values[8,46] = 0.0001146235521235521
$values[8,47] = 0.001234286783493983
$values[8,48] = 0.0008077252184532221
$values[8,49] = 0.0001146235521235521
$values[8,50] = 0.0008077252184532221
$values[8,51] = 0.0008077252184532221
$values[8,52] = 0.0001146235521235521
$values[8,53] = 0.0008077252184532221

# Row 10
$values[9,0] = 0.03239724877960055
$values[9,1] = 0.0001146235521235521
$values[9,2] = 0.03239724877960055
$values[9,3] = 0.03239724877960055
$values[9,4] = 0.0001146235521235521
$values[9,5] = 0.03239724877960055
$values[9,6] = 0.01342113595031052
$values[9,7] = 0.0001518968383129674
$values[9,8] = 0.01342113595031052
$values[9,9] = 0.02271664407576451
$values[9,10] = 0.0001146235521235521
$values[9,11] = 0.02271664407576451
$values[9,12] = 0.01342113595031052
$values[9,13] = 0.0001518968383129674
$values[9,14] = 0.01342113595031052
$values[9,15] = 0.01342113595031052
$values[9,16] = 0.0001518968383129674
$values[9,17] = 0.01342113595031052
$values[9,18] = 0.01342113595031052
$values[9,19] = 0.0001518968383129674
$values[9,20] = 0.01342113595031052
$values[9,21] = 0.01826206209161527
$values[9,22] = 0.0001394724095831623
$values[9,23] = 0.01826206209161527
$values[9,24] = 0.02271664407576451
$values[9,25] = 0.0001146235521235521
$values[9,26] = 0.02271664407576451
$values[9,27] = 0.02271664407576451
$values[9,28] = 0.0001146235521235521
$values[9,29] = 0.02271664407576451
$values[9,30] = 0.02271664407576451
$values[9,31] = 0.0001146235521235521
$values[9,32] = 0.02271664407576451
$values[9,33] = 0.01177336119396915
$values[9,34] = 0.0001146235521235521
$values[9,35] = 0.01177336119396915
$values[9,36] = 0.01177336119396915
$values[9,37] = 0.0001146235521235521
$values[9,38] = 0.01177336119396915
$values[9,39] = 0.01342113595031052
$values[9,40] = 0.0001518968383129674
$values[9,41] = 0.01342113595031052
$values[9,42] = 0.0194897758411525
$values[9,43] = 0.0001146235521235521
$values[9,44] = 0.0194897758411525
$values[9,45] = 0.0194897758411525
$values[9,46] = 0.0001146235521235521
$values[9,47] = 0.0194897758411525
$values[9,48] = 0.01177336119396915
$values[9,49] = 0.0001146235521235521
$values[9,50] = 0.01177336119396915
$values[9,51] = 0.01177336119396915
$values[9,52] = 0.0001146235521235521
$values[9,53] = 0.01177336119396915

$ws.Range("A1:BB10").Value = $values
